# Update countries & provincias Spain
# Applies the refreshed COVID-19 country data pull: updates the "last updated"
# timestamp, refreshes several countries' figures in place, and re-sorts two
# countries (Egipto, Mayotte) upward past neighbours whose totals they
# overtook (their rows get the new country + numbers, the overtaken
# countries shift down one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param([int]$Row, [int]$StartCol, [object[]]$Values)
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $ws.Cells.Item($Row, $StartCol + $i).Value = $Values[$i]
    }
}

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 20:22"

# --- Simple in-place refreshes (country stays on its row) --------------
Set-RowValues 4 2 @(804476, 11717, 76104, 684660, 14016, 1198, 43712)   # Estados Unidos
Set-RowValues 7 2 @(158050, 2667, 39181, 98073, 5433, 531, 20796)       # Francia
Set-RowValues 8 2 @(148024, 959, 95200, 47876, 2889, 86, 4948)          # Alemania
Set-RowValues 18 4 @(19400, 7185)                                      # Suiza
Set-RowValues 21 2 @(17837, 1512, 6982, 10371, 380, 39, 484)           # Peru
Set-RowValues 25 2 @(13942, 229, 4507, 9251, 139, 7, 184)              # Israel
Set-RowValues 34 4 @(839, 8275)                                        # Singapur
Set-RowValues 91 4 @(282, 398)                                         # Principado de Andorra

# --- Egipto overtakes Sudafrica y Banglades (rows 53-55) -------------
Set-RowValues 53 1 @("Egipto", 3490, 157, 870, 2356, 0, 14, 264)
Set-RowValues 54 1 @("Sudafrica", 3465, 165, 1055, 2352, 36, 0, 58)
Set-RowValues 55 1 @("Banglades", 3382, 434, 87, 3185, 1, 9, 110)

# --- Mayotte overtakes Sri Lanka, Isla de Man, Kenia, Guatemala, --------
# --- Venezuela (rows 115-120) -------------------------------------------
Set-RowValues 115 1 @("Mayotte", 311, 40, 117, 190, 4, 0, 4)
Set-RowValues 116 1 @("Sri Lanka", 310, 6, 102, 201, 1, 0, 7)
Set-RowValues 117 1 @("Isla de Man", 307, 7, 200, 98, 15, 0, 9)
Set-RowValues 118 1 @("Kenia", 296, 15, 74, 208, 2, 0, 14)
Set-RowValues 119 1 @("Guatemala", 294, 5, 24, 263, 3, 0, 7)
Set-RowValues 120 1 @("Venezuela", 285, 0, 117, 158, 4, 0, 10)
